$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Value2 -eq 45184) {
        $ws.Cells.Item($r, 3).Value = 45185
    }
}
